$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text or $null, new Volume(1h) (E) text or $null,
# and a flag indicating the Price value looks numeric and must be forced to Text so Excel
# does not silently convert it into a number (losing formatting such as trailing zeros).
$updates = @(
    @{ Row = 2; D = '38.134.35'; DForceText = $false; E = '  +0.60%  ' }
    @{ Row = 3; D = '2.090.79'; DForceText = $false; E = '  +2.69%  ' }
    @{ Row = 4; D = $null; DForceText = $false; E = '  +0.03%  ' }
    @{ Row = 5; D = '228.88'; DForceText = $true; E = '  +0.41%  ' }
    @{ Row = 6; D = '0.613'; DForceText = $true; E = '  +0.41%  ' }
    @{ Row = 7; D = '60.51'; DForceText = $true; E = '  -0.56%  ' }
    @{ Row = 8; D = $null; DForceText = $false; E = '  +0.02%  ' }
    @{ Row = 9; D = '0.378'; DForceText = $true; E = '  -0.18%  ' }
    @{ Row = 10; D = '0.0842'; DForceText = $true; E = '  +2.99%  ' }
    @{ Row = 11; D = $null; DForceText = $false; E = '  -0.23%  ' }
    @{ Row = 12; D = '2.399.76'; DForceText = $false; E = '  +2.76%  ' }
    @{ Row = 13; D = '14.59'; DForceText = $true; E = '  +0.37%  ' }
    @{ Row = 14; D = '22.18'; DForceText = $true; E = '  +3.83%  ' }
    @{ Row = 15; D = '5.47'; DForceText = $true; E = '  +5.71%  ' }
    @{ Row = 16; D = '0.773'; DForceText = $true; E = '  +1.27%  ' }
    @{ Row = 17; D = '2.073.29'; DForceText = $false; E = '  +0.89%  ' }
    @{ Row = 18; D = '38.074.50'; DForceText = $false; E = '  +0.53%  ' }
    @{ Row = 19; D = '6.00'; DForceText = $true; E = '  +1.16%  ' }
    @{ Row = 20; D = '70.09'; DForceText = $true; E = '  +0.47%  ' }
    @{ Row = 21; D = '0.0₃0833'; DForceText = $false; E = $null }
    @{ Row = 22; D = '223.92'; DForceText = $true; E = '  +0.57%  ' }
    @{ Row = 23; D = $null; DForceText = $false; E = '  -0.09%  ' }
    @{ Row = 24; D = '2.44'; DForceText = $true; E = '  +0.49%  ' }
    @{ Row = 25; D = '2.32'; DForceText = $true; E = '  +3.19%  ' }
    @{ Row = 26; D = '169.85'; DForceText = $true; E = '  +1.53%  ' }
    @{ Row = 27; D = '9.39'; DForceText = $true; E = '  +0.50%  ' }
    @{ Row = 28; D = $null; DForceText = $false; E = '  -0.01%  ' }
    @{ Row = 29; D = '18.94'; DForceText = $true; E = '  +0.27%  ' }
    @{ Row = 30; D = '1.34'; DForceText = $true; E = '  +4.86%  ' }
    @{ Row = 31; D = $null; DForceText = $false; E = '  -0.47%  ' }
    @{ Row = 32; D = $null; DForceText = $false; E = '  +4.91%  ' }
    @{ Row = 33; D = '4.69'; DForceText = $true; E = '  +4.08%  ' }
    @{ Row = 34; D = $null; DForceText = $false; E = '  +0.01%  ' }
    @{ Row = 35; D = $null; DForceText = $false; E = '  -0.49%  ' }
    @{ Row = 36; D = $null; DForceText = $false; E = '  +4.57%  ' }
    @{ Row = 37; D = '6.42'; DForceText = $true; E = '  +0.68%  ' }
    @{ Row = 38; D = '3.51'; DForceText = $true; E = '  +5.53%  ' }
    @{ Row = 39; D = '0.999'; DForceText = $true; E = '  -0.22%  ' }
    @{ Row = 40; D = '18.01'; DForceText = $true; E = '  +2.56%  ' }
    @{ Row = 41; D = '1.558.28'; DForceText = $false; E = '  +1.45%  ' }
    @{ Row = 42; D = '99.88'; DForceText = $true; E = '  +3.53%  ' }
    @{ Row = 43; D = '0.0218'; DForceText = $true; E = '  +0.18%  ' }
    @{ Row = 44; D = $null; DForceText = $false; E = '  +1.16%  ' }
    @{ Row = 45; D = '0.0914'; DForceText = $true; E = '  -0.11%  ' }
    @{ Row = 46; D = '4.15'; DForceText = $true; E = '  +3.35%  ' }
    @{ Row = 47; D = $null; DForceText = $false; E = '  +0.89%  ' }
    @{ Row = 48; D = '7.44'; DForceText = $true; E = '  +5.04%  ' }
    @{ Row = 49; D = $null; DForceText = $false; E = '  +1.21%  ' }
    @{ Row = 50; D = '2.98'; DForceText = $true; E = '  +0.72%  ' }
    @{ Row = 51; D = '2.287.78'; DForceText = $false; E = '  +2.82%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dcell = $ws.Cells.Item($u.Row, 4)
        if ($u.DForceText) {
            $dcell.NumberFormat = "@"
            $dcell.Value = $u.D
            $dcell.Style = "Normal"
        } else {
            $dcell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

Write-Host "Done updating cryptos list."
